$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.051.93"
$ws.Range("E2").Value = "  +2.39%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.673.38"
$ws.Range("E3").Value = "  +3.25%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.28"
$ws.Range("E5").Value = "  +1.58%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.511"
$ws.Range("E6").Value = "  +2.06%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.01%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +2.35%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0618"
$ws.Range("E9").Value = "  +1.50%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  +4.96%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +4.71%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.908.76"
$ws.Range("E12").Value = "  +3.29%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.675.33"
$ws.Range("E13").Value = "  +3.37%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +1.68%  "

# Row 15 - Litecoin (was Polygon)
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "65.82"
$ws.Range("E15").Value = "  +3.15%  "

# Row 16 - Polygon (was Litecoin)
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.521"
$ws.Range("E16").Value = "  +2.42%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.070.01"
$ws.Range("E17").Value = "  +2.41%  "

# Row 18 - BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.53"
$ws.Range("E18").Value = "  -0.30%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +1.63%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.72"
$ws.Range("E20").Value = "  -1.69%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.04%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.48"
$ws.Range("E22").Value = "  +3.84%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +1.81%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.17%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.34"
$ws.Range("E25").Value = "  -1.19%  "

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.16"
$ws.Range("E26").Value = "  +1.16%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  +0.56%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.91"
$ws.Range("E28").Value = "  +2.19%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.05%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +0.28%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +1.87%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +1.98%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.454.63"
$ws.Range("E33").Value = "  -4.23%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.15"
$ws.Range("E34").Value = "  +5.51%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  +6.44%  "

# Row 37 - ImmutableX
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.571"
$ws.Range("E37").Value = "  +0.60%  "

# Row 38 - ARBITRUM
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.895"

# Row 39 - VeChain
$ws.Range("E39").Value = "  +1.84%  "

# Row 40 - FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.07"
$ws.Range("E40").Value = "  +3.11%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  -0.01%  "

# Row 42 - WEMIXToken
$ws.Range("E42").Value = "  +10.58%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +3.39%  "

# Row 44 - Aave
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.21"
$ws.Range("E44").Value = "  +5.39%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.818.00"
$ws.Range("E45").Value = "  +3.27%  "

# Row 46 - TrustWalletToken
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.780"
$ws.Range("E46").Value = "  +2.30%  "

# Row 47 - Quant
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.31"
$ws.Range("E47").Value = "  -0.31%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  +1.99%  "

# Row 49 - Algorand
$ws.Range("E49").Value = "  +3.96%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +1.24%  "

# Row 51 - EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("E51").Value = "  +1.72%  "
